# Notion "Lũy kế tháng LONG XUYÊN" export update (Add thêm nhân sự Nguyễn Hữu Quang)
#
# Underlying change: two of the shared "last_edited_time" timestamp strings were
# updated (2024-08-31T05:40 -> 2024-08-31T15:46, and 2024-08-31T05:43 -> 2024-08-31T15:45),
# and two page rows (Tháng 4 / row 4, Tháng 5 / row 13) that used to reference the
# "05:43" timestamp now reference the "15:46" one instead. All cells in column D
# (last_edited_time) that share those underlying strings pick up the new text too.
# Separately, row 10 (Tháng 8) had its "Chi tiêu" (W10) and "Lũy kế" (AA10) amounts
# adjusted.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# last_edited_time (column D) updates
$ws.Range("D3").Value  = "2024-08-31T15:46:00.000Z"
$ws.Range("D4").Value  = "2024-08-31T15:46:00.000Z"
$ws.Range("D5").Value  = "2024-08-31T15:46:00.000Z"
$ws.Range("D7").Value  = "2024-08-31T15:45:00.000Z"
$ws.Range("D10").Value = "2024-08-31T15:46:00.000Z"
$ws.Range("D13").Value = "2024-08-31T15:46:00.000Z"

# Row 10 (Tháng 8) numeric updates: Chi tiêu (W10), Lũy kế (AA10)
$ws.Range("W10").Value  = 58993000
$ws.Range("AA10").Value = 61507000
